$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the four "JUAN DIEGO LOPEZ VERANO" period rows (old rows 24-27); this
# leaves the old "total" row (28) in place, which becomes row 24, and shifts
# the signature block up from rows 33/34 to rows 29/30.
$ws.Rows("24:27").Delete()

# ---- Header block updates ----
$ws.Range("E11").Value = 628100
$ws.Range("F13").Value = 5

# ---- Worker detail rows 16-21: swap which worker covers each period ----
# Row 16: CARLOS ENRIQUE FERNANDEZ ANILLO, period 2402
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1143334338"
$ws.Range("D16").Value = "CARLOS ENRIQUE FERNANDEZ ANILLO"
$ws.Range("E16").Value = "2402"
$ws.Range("F16").Value = 100000
$ws.Range("G16").Value = 2500000

# Row 17: PAULA ANDREA MATOS ROJANO, period 2402
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1235046520"
$ws.Range("D17").Value = "PAULA ANDREA MATOS ROJANO"
$ws.Range("E17").Value = "2402"
$ws.Range("F17").Value = 60000
$ws.Range("G17").Value = 1590000

# Row 18: CARLOS ENRIQUE FERNANDEZ ANILLO, period 2403
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1143334338"
$ws.Range("D18").Value = "CARLOS ENRIQUE FERNANDEZ ANILLO"
$ws.Range("E18").Value = "2403"
$ws.Range("F18").Value = 100000
$ws.Range("G18").Value = 2500000

# Row 19: PAULA ANDREA MATOS ROJANO, period 2403
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1235046520"
$ws.Range("D19").Value = "PAULA ANDREA MATOS ROJANO"
$ws.Range("E19").Value = "2403"
$ws.Range("F19").Value = 60000
$ws.Range("G19").Value = 1590000

# Row 20: CARLOS ENRIQUE FERNANDEZ ANILLO, period 2404
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1143334338"
$ws.Range("D20").Value = "CARLOS ENRIQUE FERNANDEZ ANILLO"
$ws.Range("E20").Value = "2404"
$ws.Range("F20").Value = 100000
$ws.Range("G20").Value = 2500000

# Row 21: PAULA ANDREA MATOS ROJANO, period 2404
$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1235046520"
$ws.Range("D21").Value = "PAULA ANDREA MATOS ROJANO"
$ws.Range("E21").Value = "2404"
$ws.Range("F21").Value = 60000
$ws.Range("G21").Value = 1590000

# ---- New worker rows 22-23 ----
# Row 22: SAMUEL EUGENIO PEREZ CHACON (PPT), period 2507
$ws.Range("B22").Value = "PPT"
$ws.Range("C22").Value = "1573749"
$ws.Range("D22").Value = "SAMUEL EUGENIO PEREZ CHACON"
$ws.Range("E22").Value = "2507"
$ws.Range("F22").Value = 27560
$ws.Range("G22").Value = 1590000

# Row 23: SHAILY STEFANI MONCADA VELASQUEZ (CC), period 2508
$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1063144471"
$ws.Range("D23").Value = "SHAILY STEFANI MONCADA VELASQUEZ"
$ws.Range("E23").Value = "2508"
$ws.Range("F23").Value = 56940
$ws.Range("G23").Value = 1423500

# ---- Total row (was row 28, now row 24) ----
$ws.Range("B24").Value = "PPT"
$ws.Range("C24").Value = "1573749"
$ws.Range("D24").Value = "SAMUEL EUGENIO PEREZ CHACON"
$ws.Range("E24").Value = "2508"
$ws.Range("F24").Value = 63600
$ws.Range("G24").Value = 1590000
